# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (classic blue/orange Office palette)
#                            only wired up via the notes master
#   ppt/theme/theme2.xml -> "Integral"/"Red Violet" palette
#                            wired up via the (only) slide master, so it is
#                            the theme that actually paints every slide.
#
# The authored change swaps the two palettes: the colours that used to live
# in theme1.xml now belong to theme2.xml (the design actually applied to the
# deck), and vice versa. In other words, the presentation's live/visible
# theme switches from the "Integral" (Red Violet) colours to the classic
# Office colours.
#
# PowerPoint's object model exposes the live theme's 12-slot colour scheme
# through Slide.ThemeColorScheme (it is shared by every slide because they
# all hang off the single slide master), so we push the "Office Theme"
# palette into it, one RGB swatch at a time.

function ComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) that theme1.xml currently holds.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation

# Any slide works -- ThemeColorScheme is master-level, not per-slide.
$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ComRgb $officeThemeColors[$i - 1]
}
